$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Range("D2:H2").ClearContents()
$ws.Columns.Item(5).NumberFormat = "0%"
$ws.Range("E1:E2").Clear()
